$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 236; this shifts rows 236-332 down to 237-333
$ws.Rows.Item(236).Insert()

# Populate the newly inserted row 236 with the new record's data
$ws.Cells.Item(236, 1).Value2 = 10
$ws.Cells.Item(236, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(236, 3).Value2 = "La Araucanía"
$ws.Cells.Item(236, 4).Value2 = 44755
$ws.Cells.Item(236, 5).Value2 = 9
$ws.Cells.Item(236, 6).Value2 = 100112044
$ws.Cells.Item(236, 7).Value2 = "Perejil"
$ws.Cells.Item(236, 8).Value2 = "Sin especificar"
$ws.Cells.Item(236, 9).Value2 = "Primera"
$ws.Cells.Item(236, 10).Value2 = 30
$ws.Cells.Item(236, 11).Value2 = 4600
$ws.Cells.Item(236, 12).Value2 = 4600
$ws.Cells.Item(236, 13).Value2 = 4600
$ws.Cells.Item(236, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(236, 15).Value2 = "Provincia de Cautín"
$ws.Cells.Item(236, 16).Value2 = 1533
$ws.Cells.Item(236, 17).Value2 = 3
$ws.Cells.Item(236, 18).Value2 = "Hortaliza"

Write-Output "done"
